$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Speaker column (D) values per the new abbreviated naming scheme
$ws.Range("D2").Value = "S"
$ws.Range("D3").Value = "SN"

# Every row where the speaker was "HILLARY LEWIS-WOLFSEN" becomes "T"
$teacherRows = @(4,5,6,7,10,12,13,16,17,18,19,20,21,23,24,25,26)
foreach ($r in $teacherRows) {
    $ws.Cells.Item($r, 4).Value = "T"
}
